$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=0.01368308067321777; C=0.02483744621276856; D=0.003437042236328125; E=0.01657395362854004; F=0; G=0.05996384620666504; H=0.02000150680541992; I=0.02145295143127441; J=0.01535792350769043; K=0.02404146194458008; L=0.003359413146972656; M=0.01656999588012695 }
    3 = @{ B=0.1307877540588379; C=0.04869184494018555; D=0.02769923210144043; E=0.01864128112792969; F=0.008669757843017578; G=0.01287837028503418; H=0.1359569072723389; I=0.0372157096862793; J=0.1070014953613281; K=0.03587250709533692; L=0.03759307861328125; M=0.01731629371643066 }
    4 = @{ B=0.05874981880187988; C=0.02714376449584961; D=0.03059911727905273; E=0.01584486961364746; F=0.09232025146484375; G=0.01421313285827637; H=0.0430945873260498; I=0.03295221328735352; J=0.03149204254150391; K=0.0277310848236084; L=0.06775507926940919; M=0.01273941993713379 }
    5 = @{ B=0.0382883071899414; C=0.03574857711791992; D=0.03637590408325195; E=0.03173689842224121; F=0; G=0; H=0.02858166694641113; I=0.03210873603820801; J=0.01864619255065918; K=0.02571640014648437; L=0; M=0 }
    6 = @{ B=0.6467616558074951; C=0.0276768684387207; D=0.6240752696990967; E=0.03111162185668945; F=1.294736576080322; G=0.02487802505493164; H=0.3938920497894287; I=0.0213068962097168; J=0.6603969097137451; K=0.0349393367767334; L=0.5810380935668945; M=0.01997919082641602 }
    7 = @{ B=0.9033839225769043; C=0.1154187202453613; D=0.4600280284881592; E=0.0659097671508789; F=0.6085456371307373; G=0.02822628021240235; H=0.9345259189605712; I=0.1059933662414551; J=0.3508360385894775; K=0.05184469223022461; L=0.7216342926025391; M=0.03750615119934082 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
